# Rename the "congenital" category to "misc_long_term" across all
# worksheets in the workbook (cell A3 holds the category name for
# each per-variable sheet).
$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    $cell = $ws.Range("A3")
    if ($cell.Value2 -eq "congenital") {
        $cell.Value = "misc_long_term"
    }
}
